$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column G, matching width of the other "month" columns' family but wider (17 chars)
$ws.Columns.Item(7).ColumnWidth = 113/7

# G1: header label "PRESUPUESTO" - copy F1's format (bold header style) then set its text
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "PRESUPUESTO"

# G2: numeric 0, reuse F2's currency-style format
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = 0

# G3: numeric 0, reuse F3's currency-style format
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = 0

$excel.CutCopyMode = $false
